## Generate Report for Handback
## - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
##   everywhere that status text is shown (Overview + per-locale sheets).
## - Each locale sheet (zh-cn, de-de) gets its "Latest Target File" (F) and
##   "Latest Handback File" (G) columns populated with hyperlinked file
##   names (mirroring the existing Source File Name / Latest Handoff File
##   hyperlinks), and the "Latest Handback DateTime" (H) column stamped
##   with the real handback timestamp instead of the zero date.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdUrlPrefix  = "https://github.com/OpenLocalizationTest/oltest/blob/7a685d1e0d486ebaa064b8eb8d2d1d344912b74f/e2e/"
$zhXlfUrlPrefix = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56153545c5a5467726974a5ed2466ef407aff28b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deXlfUrlPrefix = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4f10b91c9817cd5196101fffdd4f4ce80b00810/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$file1Md  = "89f49ef1-90e6-4b57-b259-76d536eb46e4.md"
$file1ZhXlf = "89f49ef1-90e6-4b57-b259-76d536eb46e4.b81e9dd5ca5fae35ab98023cec38426ecbd87773.zh-cn.xlf"
$file1DeXlf = "89f49ef1-90e6-4b57-b259-76d536eb46e4.b81e9dd5ca5fae35ab98023cec38426ecbd87773.de-de.xlf"

$file2Md  = "faa2eff6-2430-4173-810c-d843ff92cbe2.md"
$file2ZhXlf = "faa2eff6-2430-4173-810c-d843ff92cbe2.2780d3b7e7070e816344dc4b347571ea7f514742.zh-cn.xlf"
$file2DeXlf = "faa2eff6-2430-4173-810c-d843ff92cbe2.2780d3b7e7070e816344dc4b347571ea7f514742.de-de.xlf"

# ---- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ----

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# ---- Helper-less inline logic: fill Latest Target File / Latest Handback
#      File / Latest Handback DateTime for each locale sheet ----

function Add-HandbackHyperlink($ws, $cellRef, $displayText, $url) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = 15570276
}

# --- zh-cn sheet ---

$url = $mdUrlPrefix + $file1Md
Add-HandbackHyperlink $wsZh "F2" $file1Md $url
$url = $zhXlfUrlPrefix + $file1ZhXlf
Add-HandbackHyperlink $wsZh "G2" $file1ZhXlf $url
$wsZh.Range("H2").Value = "2016-03-22 16:52:28"

$url = $mdUrlPrefix + $file2Md
Add-HandbackHyperlink $wsZh "F3" $file2Md $url
$url = $zhXlfUrlPrefix + $file2ZhXlf
Add-HandbackHyperlink $wsZh "G3" $file2ZhXlf $url
$wsZh.Range("H3").Value = "2016-03-22 16:52:28"

# --- de-de sheet ---

$url = $mdUrlPrefix + $file1Md
Add-HandbackHyperlink $wsDe "F2" $file1Md $url
$url = $deXlfUrlPrefix + $file1DeXlf
Add-HandbackHyperlink $wsDe "G2" $file1DeXlf $url
$wsDe.Range("H2").Value = "2016-03-22 16:52:37"

$url = $mdUrlPrefix + $file2Md
Add-HandbackHyperlink $wsDe "F3" $file2Md $url
$url = $deXlfUrlPrefix + $file2DeXlf
Add-HandbackHyperlink $wsDe "G3" $file2DeXlf $url
$wsDe.Range("H3").Value = "2016-03-22 16:52:37"
